# Auto-generated edit script: update crypto price/volume table
# per commit "Updated cryptos list on Wed Aug  2 22:17:35 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.182.58'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.842.94'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'241.88"
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("D6").Value = "'0.6880"
$ws.Range("E6").Value = '  -2.06%  '
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = "'0.3019"
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("D9").Value = "'0.07488"
$ws.Range("E9").Value = '  -3.50%  '
$ws.Range("D10").Value = "'23.24"
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").Value = "'0.07662"
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").Value = '1.842.32'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = "'5.070"
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("D14").Value = "'0.6849"
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = "'87.63"
$ws.Range("E15").Value = '  -5.95%  '
$ws.Range("D16").Value = "'6.191"
$ws.Range("E16").Value = '  -6.17%  '
$ws.Range("D17").Value = '29.176.37'
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").Value = "'0.000008192"
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").Value = '2.081.77'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = "'230.09"
$ws.Range("E20").Value = '  -4.91%  '
$ws.Range("D21").Value = "'12.57"
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = "'7.421"
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("D24").Value = "'0.9991"
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  -3.76%  '
$ws.Range("D26").Value = "'159.35"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = "'8.794"
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("D28").Value = "'18.14"
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").Value = "'1.518"
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = "'4.291"
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("D31").Value = "'4.151"
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("D32").Value = "'1.196"
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = "'0.05241"
$ws.Range("E33").Value = '  +2.30%  '
$ws.Range("D34").Value = "'0.7617"
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("D35").Value = "'1.859"
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("D36").Value = "'1.137"
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").Value = '1.306.64'
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("D39").Value = "'0.01837"
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("D40").Value = "'2.727"
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").Value = "'0.9361"
$ws.Range("E41").Value = '  -1.89%  '
$ws.Range("D42").Value = "'5.968"
$ws.Range("E42").Value = '  -1.85%  '
$ws.Range("D43").Value = "'105.18"
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("D44").Value = "'0.9987"
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.985.78'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'65.19"
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("D47").Value = "'0.5196"
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'9.560"
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = "'0.00000000122"
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("D50").Value = "'1.779"
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("D51").Value = "'0.05955"
$ws.Range("E51").Value = '  +0.82%  '
